# Refactor synthetic array /3: rename the "black/red/green/orange" emoji+label
# columns to "blue/red/green/orange" (bleu instead of noir), updating both the
# pictogram column (A) and the color-name column (B) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old -> new values for the pictogram column (A)
$pictoMap = @{
    "⬛" = "📘"
    "🟥" = "📕"
    "🟩" = "📗"
    "🟧" = "📙"
}

# Mapping of old -> new values for the label column (B); only "noir" changes
$labelMap = @{
    "noir" = "bleu"
}

# Determine the extent of the used range
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value2
    if ($null -ne $valA -and $pictoMap.ContainsKey($valA)) {
        $cellA.Value2 = $pictoMap[$valA]
    }

    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value2
    if ($null -ne $valB -and $labelMap.ContainsKey($valB)) {
        $cellB.Value2 = $labelMap[$valB]
    }
}
